$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 7) with the same data pattern as row 6,
# changing the Total Price value to 312.
$ws.Range("A7").Value = "Elmar Qara"
$ws.Range("B7").Value = "elmarqarayev69@gmail.com"
$ws.Range("C7").Value = 312
$ws.Range("D7").Value = "Pending"
